# "added harvard case classification"
#
# The underlying dataset gained a new case-classification dimension, which
# shifts the previously-rightmost summary column ("average_doctor") one slot
# to the right (old BP -> BQ) and introduces a new "average_doctor_old"
# summary column in its place (new BP), and recomputes every per-app
# average/variance/std-dev statistic against the updated classification.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: average_doctor / average_doctor_old swap columns BP/BQ ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Recomputed per-row statistics (columns affected by the new case classification) ---
# Row 4
$ws.Range("E4").Value = 0.42
$ws.Range("F4").Value = 0.07199999999999999
$ws.Range("G4").Value = 0.268
$ws.Range("N4").Value = 0.433
$ws.Range("O4").Value = 0.065
$ws.Range("P4").Value = 0.255
$ws.Range("Q4").Value = 0.024
$ws.Range("R4").Value = 0.017
$ws.Range("S4").Value = 0.129
$ws.Range("W4").Value = 0.286
$ws.Range("X4").Value = 0.109
$ws.Range("Y4").Value = 0.33
$ws.Range("AI4").Value = 0.292
$ws.Range("AJ4").Value = 0.08799999999999999
$ws.Range("AK4").Value = 0.297
$ws.Range("AU4").Value = 0.19
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.167
$ws.Range("BA4").Value = 2.007
$ws.Range("BB4").Value = 0.16
$ws.Range("BC4").Value = 0.4
$ws.Range("BG4").Value = 0.729
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.376
$ws.Range("BM4").Value = 0.716
$ws.Range("BN4").Value = 0.08
$ws.Range("BO4").Value = 0.282
$ws.Range("BP4").Value = 0.669
$ws.Range("BQ4").Value = 0.709
# Row 5
$ws.Range("E5").Value = 0.537
$ws.Range("F5").Value = 0.08599999999999999
$ws.Range("G5").Value = 0.294
$ws.Range("N5").Value = 0.733
$ws.Range("O5").Value = 0.082
$ws.Range("P5").Value = 0.286
$ws.Range("Q5").Value = 0.016
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.083
$ws.Range("W5").Value = 0.276
$ws.Range("X5").Value = 0.109
$ws.Range("Y5").Value = 0.33
$ws.Range("AI5").Value = 0.312
$ws.Range("AJ5").Value = 0.099
$ws.Range("AK5").Value = 0.314
$ws.Range("AU5").Value = 0.369
$ws.Range("AV5").Value = 0.093
$ws.Range("AW5").Value = 0.305
$ws.Range("BA5").Value = 1.337
$ws.Range("BC5").Value = 0.282
$ws.Range("BG5").Value = 0.395
$ws.Range("BI5").Value = 0.225
$ws.Range("BM5").Value = 0.551
$ws.Range("BN5").Value = 0.064
$ws.Range("BO5").Value = 0.252
$ws.Range("BP5").Value = 0.446
$ws.Range("BQ5").Value = 0.456
# Row 6
$ws.Range("E6").Value = 0.471
$ws.Range("N6").Value = 0.544
$ws.Range("Q6").Value = 0.019
$ws.Range("W6").Value = 0.281
$ws.Range("AI6").Value = 0.302
$ws.Range("AU6").Value = 0.251
$ws.Range("BA6").Value = 1.596
$ws.Range("BG6").Value = 0.512
$ws.Range("BM6").Value = 0.623
$ws.Range("BP6").Value = 0.532
$ws.Range("BQ6").Value = 0.552
# Row 7
$ws.Range("E7").Value = 0.509
$ws.Range("N7").Value = 0.644
$ws.Range("Q7").Value = 0.017
$ws.Range("W7").Value = 0.278
$ws.Range("AI7").Value = 0.308
$ws.Range("AU7").Value = 0.31
$ws.Range("BA7").Value = 1.429
$ws.Range("BG7").Value = 0.435
$ws.Range("BM7").Value = 0.578
$ws.Range("BP7").Value = 0.476
$ws.Range("BQ7").Value = 0.49
# Row 8
$ws.Range("E8").Value = 0.602
$ws.Range("F8").Value = 0.112
$ws.Range("G8").Value = 0.335
$ws.Range("N8").Value = 0.773
$ws.Range("O8").Value = 0.066
$ws.Range("P8").Value = 0.258
$ws.Range("Q8").Value = 0.017
$ws.Range("W8").Value = 0.304
$ws.Range("AI8").Value = 0.334
$ws.Range("AJ8").Value = 0.129
$ws.Range("AK8").Value = 0.359
$ws.Range("AU8").Value = 0.311
$ws.Range("AV8").Value = 0.08400000000000001
$ws.Range("AW8").Value = 0.291
$ws.Range("BA8").Value = 1.748
$ws.Range("BB8").Value = 0.125
$ws.Range("BC8").Value = 0.353
$ws.Range("BG8").Value = 0.5649999999999999
$ws.Range("BH8").Value = 0.107
$ws.Range("BI8").Value = 0.326
$ws.Range("BM8").Value = 0.695
$ws.Range("BN8").Value = 0.066
$ws.Range("BO8").Value = 0.256
$ws.Range("BP8").Value = 0.583
$ws.Range("BQ8").Value = 0.605
# Row 9
$ws.Range("E9").Value = 0.548
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.677
$ws.Range("O9").Value = 0.219
$ws.Range("P9").Value = 0.467
$ws.Range("W9").Value = 0.204
$ws.Range("X9").Value = 0.163
$ws.Range("Y9").Value = 0.403
$ws.Range("AI9").Value = 0.258
$ws.Range("AJ9").Value = 0.191
$ws.Range("AK9").Value = 0.438
$ws.Range("BA9").Value = 1.71
$ws.Range("BB9").Value = 0.248
$ws.Range("BC9").Value = 0.498
$ws.Range("BG9").Value = 0.602
$ws.Range("BH9").Value = 0.24
$ws.Range("BI9").Value = 0.489
$ws.Range("BM9").Value = 0.656
$ws.Range("BN9").Value = 0.226
$ws.Range("BO9").Value = 0.475
$ws.Range("BP9").Value = 0.57
$ws.Range("BQ9").Value = 0.588
# Row 10
$ws.Range("E10").Value = 0.677
$ws.Range("F10").Value = 0.219
$ws.Range("G10").Value = 0.467
$ws.Range("N10").Value = 0.871
$ws.Range("O10").Value = 0.112
$ws.Range("P10").Value = 0.335
$ws.Range("W10").Value = 0.376
$ws.Range("X10").Value = 0.235
$ws.Range("Y10").Value = 0.484
$ws.Range("AI10").Value = 0.366
$ws.Range("AJ10").Value = 0.232
$ws.Range("AK10").Value = 0.482
$ws.Range("AU10").Value = 0.301
$ws.Range("AV10").Value = 0.21
$ws.Range("AW10").Value = 0.459
$ws.Range("BA10").Value = 2.086
$ws.Range("BB10").Value = 0.243
$ws.Range("BC10").Value = 0.493
$ws.Range("BG10").Value = 0.656
$ws.Range("BH10").Value = 0.226
$ws.Range("BI10").Value = 0.475
$ws.Range("BM10").Value = 0.849
$ws.Range("BN10").Value = 0.128
$ws.Range("BO10").Value = 0.358
$ws.Range("BP10").Value = 0.695
$ws.Range("BQ10").Value = 0.726
# Row 11
$ws.Range("E11").Value = 0.71
$ws.Range("F11").Value = 0.206
$ws.Range("G11").Value = 0.454
$ws.Range("N11").Value = 0.892
$ws.Range("O11").Value = 0.096
$ws.Range("P11").Value = 0.31
$ws.Range("W11").Value = 0.376
$ws.Range("X11").Value = 0.235
$ws.Range("Y11").Value = 0.484
$ws.Range("AI11").Value = 0.398
$ws.Range("AJ11").Value = 0.24
$ws.Range("AK11").Value = 0.489
$ws.Range("AU11").Value = 0.441
$ws.Range("AV11").Value = 0.247
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.086
$ws.Range("BB11").Value = 0.243
$ws.Range("BC11").Value = 0.493
$ws.Range("BG11").Value = 0.656
$ws.Range("BH11").Value = 0.226
$ws.Range("BI11").Value = 0.475
$ws.Range("BM11").Value = 0.849
$ws.Range("BN11").Value = 0.128
$ws.Range("BO11").Value = 0.358
$ws.Range("BP11").Value = 0.695
$ws.Range("BQ11").Value = 0.729
# Row 12
$ws.Range("E12").Value = 1.409
$ws.Range("F12").Value = 0.757
$ws.Range("G12").Value = 0.87
$ws.Range("N12").Value = 1.471
$ws.Range("O12").Value = 1.049
$ws.Range("P12").Value = 1.024
$ws.Range("W12").Value = 1.629
$ws.Range("X12").Value = 0.576
$ws.Range("Y12").Value = 0.759
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.767
$ws.Range("AV12").Value = 2.737
$ws.Range("AW12").Value = 1.654
$ws.Range("BA12").Value = 3.704
$ws.Range("BB12").Value = 0.401
$ws.Range("BC12").Value = 0.633
$ws.Range("BG12").Value = 1.098
$ws.Range("BH12").Value = 0.121
$ws.Range("BI12").Value = 0.349
$ws.Range("BM12").Value = 1.291
$ws.Range("BN12").Value = 0.333
$ws.Range("BO12").Value = 0.577
$ws.Range("BP12").Value = 1.235
$ws.Range("BQ12").Value = 1.26
# Row 13
$ws.Range("E13").Value = 1.579
$ws.Range("F13").Value = 0.656
$ws.Range("G13").Value = 0.8100000000000001
$ws.Range("N13").Value = 2.069
$ws.Range("O13").Value = 0.9379999999999999
$ws.Range("P13").Value = 0.969
$ws.Range("W13").Value = 1.037
$ws.Range("X13").Value = 0.193
$ws.Range("Y13").Value = 0.439
$ws.Range("AI13").Value = 1.28
$ws.Range("AJ13").Value = 0.37
$ws.Range("AK13").Value = 0.608
$ws.Range("AU13").Value = 2.285
$ws.Range("AV13").Value = 0.925
$ws.Range("AW13").Value = 0.962
$ws.Range("BA13").Value = 2.352
$ws.Range("BB13").Value = 0.298
$ws.Range("BC13").Value = 0.546
$ws.Range("BG13").Value = 0.583
$ws.Range("BH13").Value = 0.07099999999999999
$ws.Range("BI13").Value = 0.267
$ws.Range("BM13").Value = 0.898
$ws.Range("BN13").Value = 0.281
$ws.Range("BO13").Value = 0.53
$ws.Range("BP13").Value = 0.784
$ws.Range("BQ13").Value = 0.726
